$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 176.28572
$ws.Range("I8").Value = 39
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 117
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 22
$ws.Range("N8").Value = -3278
$ws.Range("H86").Value = 7598.2383
$ws.Range("I86").Value = 2645.6365
$ws.Range("K86").Value = 2645.6365
$ws.Range("M86").Value = -1522.6365
$ws.Range("H89").Value = 7598.2383
$ws.Range("I89").Value = 2645.6365
$ws.Range("K89").Value = 13228.1825
$ws.Range("M89").Value = -7612.182500000001
$ws.Range("H111").Value = 2965.4285
$ws.Range("I111").Value = 3791.6
$ws.Range("J111").Value = 900
$ws.Range("K111").Value = 11374.8
$ws.Range("L111").Value = 2700
$ws.Range("M111").Value = -8307.799999999999
$ws.Range("N111").Value = -8834
$ws.Range("H113").Value = 52635772
$ws.Range("I113").Value = 100003880
$ws.Range("J113").Value = 4537.5557
$ws.Range("K113").Value = 100003880
$ws.Range("L113").Value = 4537.5557
$ws.Range("M113").Value = -100000626
$ws.Range("N113").Value = -11045.5557
$ws.Range("H141").Value = 2082.5
$ws.Range("I141").Value = 1853.125
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 5559.375
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -379.375
$ws.Range("N141").Value = -19360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1592.84
$ws.Range("I2").Value = 1677.4286
$ws.Range("K2").Value = 1677.4286
$ws.Range("M2").Value = -1564.4286
$ws.Range("H45").Value = 2452.261
$ws.Range("I45").Value = 3323.8
$ws.Range("J45").Value = 1781.8462
$ws.Range("K45").Value = 3323.8
$ws.Range("L45").Value = 1781.8462
$ws.Range("M45").Value = -2946.8
$ws.Range("N45").Value = -2535.8462
$ws.Range("H116").Value = 1592.84
$ws.Range("I116").Value = 1677.4286
$ws.Range("K116").Value = 1677.4286
$ws.Range("M116").Value = 616.5714
$ws.Range("H132").Value = 14342.75
$ws.Range("I132").Value = 1842.7858
$ws.Range("K132").Value = 5528.357400000001
$ws.Range("M132").Value = -2998.357400000001
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1592.84
$ws.Range("I3").Value = 1677.4286
$ws.Range("K3").Value = 1677.4286
$ws.Range("M3").Value = -1563.4286
$ws.Range("H94").Value = 2438.6667
$ws.Range("I94").Value = 2276.5881
$ws.Range("K94").Value = 2276.5881
$ws.Range("M94").Value = -1825.5881
$ws.Range("H99").Value = 1589.0625
$ws.Range("I99").Value = 1042.5
$ws.Range("K99").Value = 1042.5
$ws.Range("M99").Value = 455.5
$ws.Range("H134").Value = 3568.394
$ws.Range("I134").Value = 3704.1614
$ws.Range("K134").Value = 11112.4842
$ws.Range("M134").Value = -8577.484199999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12711.146
$ws.Range("I31").Value = 18165.959
$ws.Range("J31").Value = 5010.2354
$ws.Range("K31").Value = 18165.959
$ws.Range("L31").Value = 5010.2354
$ws.Range("M31").Value = -17870.959
$ws.Range("N31").Value = -5600.2354
$ws.Range("H34").Value = 12711.146
$ws.Range("I34").Value = 18165.959
$ws.Range("J34").Value = 5010.2354
$ws.Range("K34").Value = 18165.959
$ws.Range("L34").Value = 5010.2354
$ws.Range("M34").Value = -17963.959
$ws.Range("N34").Value = -5414.2354
$ws.Range("H58").Value = 20553.73
$ws.Range("I58").Value = 1458
$ws.Range("J58").Value = 56623.445
$ws.Range("K58").Value = 1458
$ws.Range("L58").Value = 56623.445
$ws.Range("M58").Value = -1255
$ws.Range("N58").Value = -57029.445
$ws.Range("H105").Value = 1946.125
$ws.Range("I105").Value = 1928.1666
$ws.Range("K105").Value = 1928.1666
$ws.Range("M105").Value = -181.1666
$ws.Range("H134").Value = 1205.4359
$ws.Range("I134").Value = 909.3125
$ws.Range("K134").Value = 2727.9375
$ws.Range("M134").Value = -192.9375
$ws.Range("H136").Value = 20553.73
$ws.Range("I136").Value = 1458
$ws.Range("J136").Value = 56623.445
$ws.Range("K136").Value = 4374
$ws.Range("L136").Value = 169870.335
$ws.Range("M136").Value = -1824
$ws.Range("N136").Value = -174970.335

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4361
$ws.Range("J63").Value = 5505.6
$ws.Range("L63").Value = 16516.8
$ws.Range("N63").Value = -18014.8
$ws.Range("H66").Value = 4361
$ws.Range("J66").Value = 5505.6
$ws.Range("L66").Value = 49550.4
$ws.Range("N66").Value = -57038.4
$ws.Range("H69").Value = 2500
$ws.Range("J69").Value = 2500
$ws.Range("L69").Value = 7500
$ws.Range("N69").Value = -9122
$ws.Range("H72").Value = 2500
$ws.Range("J72").Value = 2500
$ws.Range("L72").Value = 22500
$ws.Range("N72").Value = -30612
$ws.Range("H125").Value = 4571.4
$ws.Range("J125").Value = 4571.4
$ws.Range("L125").Value = 13714.2
$ws.Range("N125").Value = -23554.2
$ws.Range("H131").Value = 778.9794000000001
$ws.Range("J131").Value = 778.9794000000001
$ws.Range("L131").Value = 2336.9382
$ws.Range("N131").Value = -12416.9382
$ws.Range("H137").Value = 27781612
$ws.Range("J137").Value = 37041824
$ws.Range("L137").Value = 111125472
$ws.Range("N137").Value = -111135672
$ws.Range("H139").Value = 3022.75
$ws.Range("I139").Value = 1482.8462
$ws.Range("K139").Value = 4448.5386
$ws.Range("M139").Value = 691.4614000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 4273689.5
$ws.Range("J107").Value = 25641224
$ws.Range("L107").Value = 25641224
$ws.Range("N107").Value = -25645064
$ws.Range("H117").Value = 20000
$ws.Range("J117").Value = 20000
$ws.Range("L117").Value = 20000
$ws.Range("N117").Value = -26884
$ws.Range("H122").Value = 60607284
$ws.Range("I122").Value = 19609126
$ws.Range("J122").Value = 200001020
$ws.Range("K122").Value = 58827378
$ws.Range("L122").Value = 600003060
$ws.Range("M122").Value = -58824928
$ws.Range("N122").Value = -600007960
$ws.Range("H132").Value = 42279.355
$ws.Range("I132").Value = 7879.222
$ws.Range("J132").Value = 104199.6
$ws.Range("K132").Value = 23637.666
$ws.Range("L132").Value = 312598.8
$ws.Range("M132").Value = -21107.666
$ws.Range("N132").Value = -317658.8
$ws.Range("H135").Value = 49960
$ws.Range("J135").Value = 49960
$ws.Range("L135").Value = 49960
$ws.Range("N135").Value = -60100

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2102.1667
$ws.Range("I93").Value = 2152.7856
$ws.Range("J93").Value = 1925
$ws.Range("K93").Value = 2152.7856
$ws.Range("L93").Value = 1925
$ws.Range("M93").Value = -904.7856000000002
$ws.Range("N93").Value = -4421
$ws.Range("H122").Value = 1403358.9
$ws.Range("I122").Value = 1636445.4
$ws.Range("J122").Value = 4840
$ws.Range("K122").Value = 4909336.199999999
$ws.Range("L122").Value = 14520
$ws.Range("M122").Value = -4906886.199999999
$ws.Range("N122").Value = -19420

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 39945
$ws.Range("J116").Value = 39945
$ws.Range("L116").Value = 39945
$ws.Range("N116").Value = -49123
$ws.Range("H122").Value = 1964.3529
$ws.Range("I122").Value = 1976.6666
$ws.Range("J122").Value = 1872
$ws.Range("K122").Value = 5929.9998
$ws.Range("L122").Value = 5616
$ws.Range("M122").Value = -3479.9998
$ws.Range("N122").Value = -10516
$ws.Range("H126").Value = 1067.8572
$ws.Range("I126").Value = 887.5
$ws.Range("J126").Value = 1140
$ws.Range("K126").Value = 2662.5
$ws.Range("L126").Value = 3420
$ws.Range("M126").Value = -192.5
$ws.Range("N126").Value = -8360
$ws.Range("H136").Value = 40002332
$ws.Range("I136").Value = 62502030
$ws.Range("J136").Value = 2867.7778
$ws.Range("K136").Value = 187506090
$ws.Range("L136").Value = 8603.3334
$ws.Range("M136").Value = -187503540
$ws.Range("N136").Value = -13703.3334
